$d = $word.ActiveDocument

# 1) caseQuery line number 587 -> 591
$old1 = "M2DocEvaluator.caseQuery(M2DocEvaluator.java:587)"
$new1 = "M2DocEvaluator.caseQuery(M2DocEvaluator.java:591)"
$r1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $r1) { Write-Output "NOT FOUND 1" }

# 2) doSwitch(1242) + caseBlock(1467) -> doSwitch(1331) + caseBlock(1556)
$old2 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1467)"
$new2 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1556)"
$r2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
if (-not $r2) { Write-Output "NOT FOUND 2" }

# 3) doSwitch(1242) + caseDocumentTemplate(297) -> doSwitch(1331) + caseDocumentTemplate(301)
$old3 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:297)"
$new3 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:301)"
$r3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
if (-not $r3) { Write-Output "NOT FOUND 3" }

# 4) doSwitch(1242)+generate(282)+M2DocUtils.generate(845)+prepareoutputAndGenerate(514)+generation(421)+GeneratedMethodAccessor73
#    -> doSwitch(1331)+generate(286)+M2DocUtils.generate(853)+prepareoutputAndGenerate(518)+generation(414)+GeneratedMethodAccessor5
$old4 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1242)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:282)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:845)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:514)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:421)`n`tat sun.reflect.GeneratedMethodAccessor73.invoke(Unknown Source)"
$new4 = "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1331)`n`tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:286)`n`tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:853)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:518)`n`tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:414)`n`tat sun.reflect.GeneratedMethodAccessor5.invoke(Unknown Source)"
$r4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
if (-not $r4) { Write-Output "NOT FOUND 4" }

# 5) Replace the eclipse JDT test runner stack trace tail with the maven surefire / tycho / equinox tail
$old5 = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"
$new5 = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:365)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeWithRerun(JUnit4Provider.java:273)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:238)`n`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:159)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:206)`n`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:161)`n`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:84)`n`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:113)`n`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:593)`n`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:205)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:137)`n`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:107)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:401)`n`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:255)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n`tat java.lang.reflect.Method.invoke(Method.java:498)`n`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:657)`n`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:594)`n`tat org.eclipse.equinox.launcher.Main.run(Main.java:1447)`n`tat org.eclipse.equinox.launcher.Main.main(Main.java:1420)"
$r5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
if (-not $r5) { Write-Output "NOT FOUND 5" }

Write-Output "Done"
